# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new sheet "2022-Q1" right before "总计", copying the
#    header formatting (style) from the "2021-Q4" sheet so it matches
#    the existing fund-holdings table layout.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Copy header row (B1:H1) including formatting from 2021-Q4.
$q4Sheet.Range("B1:H1").Copy($newSheet.Range("B1"))

# Fund holding rows for 2022-Q1.
$fundRows = @(
    @{ idx = 0; code = "516110"; name = "国泰中证800汽车与零部件ETF"; scale = "1.74"; pos = "97.85"; ratio = "3.07"; value = "0.0534"; rank = 9 },
    @{ idx = 1; code = "007713"; name = "华富科技动能混合";           scale = "0.56"; pos = "86.98"; ratio = "4.84"; value = "0.0271"; rank = 7 },
    @{ idx = 2; code = "010711"; name = "华富国潮优选混合";           scale = "0.25"; pos = "94.50"; ratio = "4.66"; value = "0.0116"; rank = 5 }
)

$row = 2
foreach ($f in $fundRows) {
    # Column A carries the same style (bold, centered, bordered) as
    # the rest of the "index" column in sibling sheets - copy it from
    # 2021-Q4!A2 so the style index matches exactly.
    $q4Sheet.Range("A2").Copy($newSheet.Cells.Item($row, 1))
    $newSheet.Cells.Item($row, 1).Value = $f.idx

    # B..G are numeric-looking but stored as plain text in the source
    # data (fund codes must keep leading zeros, percentages keep their
    # literal formatting). Force text with a leading apostrophe, then
    # strip the auto-applied "quote prefix" styling by pasting the
    # (unstyled) format from the matching 2021-Q4 data row on top.
    $newSheet.Cells.Item($row, 2).Value = "'" + $f.code
    $newSheet.Cells.Item($row, 3).Value = $f.name
    $newSheet.Cells.Item($row, 4).Value = "'" + $f.scale
    $newSheet.Cells.Item($row, 5).Value = "'" + $f.pos
    $newSheet.Cells.Item($row, 6).Value = "'" + $f.ratio
    $newSheet.Cells.Item($row, 7).Value = "'" + $f.value
    $newSheet.Cells.Item($row, 8).Value = $f.rank

    $q4Sheet.Range("B2:H2").Copy()
    $newSheet.Range($newSheet.Cells.Item($row, 2), $newSheet.Cells.Item($row, 8)).PasteSpecial(-4122)

    $row = $row + 1
}

# ------------------------------------------------------------------
# 2. Insert a new summary row into "总计" for 2022-Q1 (3 funds held,
#    0.09 billion yuan), pushing the existing rows down by one.
#    NOTE: re-fetch the "总计" worksheet by name - the sheet position
#    shifted once the new "2022-Q1" sheet was inserted in front of it,
#    so the earlier $totalSheet handle no longer points at it.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Row-insert carries the formatting of the row below into the new
# blank row; reset B2:D2 back to the plain (unstyled) look used by the
# rest of the data rows, and give A2 the same bold/centered/bordered
# style used by every other cell in column A.
$totalSheet.Range("B3:D3").Copy($totalSheet.Range("B2"))
$totalSheet.Range("A3").Copy($totalSheet.Range("A2"))
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 3
$totalSheet.Cells.Item(2, 4).Value = 0.09

# Renumber the rows pushed down beneath the new one (old 0..4 -> 1..5).
for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

Write-Host "done"
